# Update the "time_taken" timestamps (column F) on the "data" sheet
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F2").Value = "2021-10-05 14:19:24.019589"
$ws1.Range("F3").Value = "2021-10-05 14:19:24.019596"
$ws1.Range("F4").Value = "2021-10-05 14:19:24.019599"
$ws1.Range("F5").Value = "2021-10-05 14:19:24.019602"
$ws1.Range("F6").Value = "2021-10-05 14:19:24.019605"
$ws1.Range("F7").Value = "2021-10-05 14:19:24.019608"
$ws1.Range("F8").Value = "2021-10-05 14:19:24.019610"
$ws1.Range("F9").Value = "2021-10-05 14:19:24.019613"
$ws1.Range("F10").Value = "2021-10-05 14:19:24.019615"
$ws1.Range("F11").Value = "2021-10-05 14:19:24.019618"
$ws1.Range("F12").Value = "2021-10-05 14:19:24.019620"
$ws1.Range("F13").Value = "2021-10-05 14:19:24.019623"
$ws1.Range("F14").Value = "2021-10-05 14:19:24.019625"
$ws1.Range("F15").Value = "2021-10-05 14:19:24.019628"
$ws1.Range("F16").Value = "2021-10-05 14:19:24.019630"
$ws1.Range("F17").Value = "2021-10-05 14:19:24.019633"
$ws1.Range("F18").Value = "2021-10-05 14:19:24.019636"
$ws1.Range("F19").Value = "2021-10-05 14:19:24.019638"
$ws1.Range("F20").Value = "2021-10-05 14:19:24.019641"
$ws1.Range("F21").Value = "2021-10-05 14:19:24.019643"
$ws1.Range("F22").Value = "2021-10-05 14:19:24.019646"
$ws1.Range("F23").Value = "2021-10-05 14:19:24.019648"
$ws1.Range("F24").Value = "2021-10-05 14:19:24.019651"
$ws1.Range("F25").Value = "2021-10-05 14:19:24.019653"
$ws1.Range("F26").Value = "2021-10-05 14:19:24.019656"
$ws1.Range("F27").Value = "2021-10-05 14:19:24.019659"
$ws1.Range("F28").Value = "2021-10-05 14:19:24.019661"
$ws1.Range("F29").Value = "2021-10-05 14:19:24.019664"
$ws1.Range("F30").Value = "2021-10-05 14:19:24.019666"
$ws1.Range("F31").Value = "2021-10-05 14:19:24.019669"
$ws1.Range("F32").Value = "2021-10-05 14:19:24.019671"
$ws1.Range("F33").Value = "2021-10-05 14:19:24.019674"
$ws1.Range("F34").Value = "2021-10-05 14:19:24.019677"
$ws1.Range("F35").Value = "2021-10-05 14:19:24.019679"
$ws1.Range("F36").Value = "2021-10-05 14:19:24.019682"
$ws1.Range("F37").Value = "2021-10-05 14:19:24.019684"
$ws1.Range("F38").Value = "2021-10-05 14:19:24.019686"
$ws1.Range("F39").Value = "2021-10-05 14:19:24.019689"
$ws1.Range("F40").Value = "2021-10-05 14:19:24.019691"
$ws1.Range("F41").Value = "2021-10-05 14:19:24.019694"
$ws1.Range("F42").Value = "2021-10-05 14:19:24.019697"
$ws1.Range("F43").Value = "2021-10-05 14:19:24.019699"
$ws1.Range("F44").Value = "2021-10-05 14:19:24.019702"
$ws1.Range("F45").Value = "2021-10-05 14:19:24.019704"
$ws1.Range("F46").Value = "2021-10-05 14:19:24.019707"
$ws1.Range("F47").Value = "2021-10-05 14:19:24.019709"
$ws1.Range("F48").Value = "2021-10-05 14:19:24.019712"
$ws1.Range("F49").Value = "2021-10-05 14:19:24.019714"
$ws1.Range("F50").Value = "2021-10-05 14:19:24.019716"
$ws1.Range("F51").Value = "2021-10-05 14:19:24.019719"
$ws1.Range("F52").Value = "2021-10-05 14:19:24.019721"
$ws1.Range("F53").Value = "2021-10-05 14:19:24.019724"
$ws1.Range("F54").Value = "2021-10-05 14:19:24.019727"
$ws1.Range("F55").Value = "2021-10-05 14:19:24.019729"
$ws1.Range("F56").Value = "2021-10-05 14:19:24.019732"
$ws1.Range("F57").Value = "2021-10-05 14:19:24.019734"
$ws1.Range("F58").Value = "2021-10-05 14:19:24.019737"
$ws1.Range("F59").Value = "2021-10-05 14:19:24.019739"
$ws1.Range("F60").Value = "2021-10-05 14:19:24.019742"
$ws1.Range("F61").Value = "2021-10-05 14:19:24.019744"
$ws1.Range("F62").Value = "2021-10-05 14:19:24.019747"
$ws1.Range("F63").Value = "2021-10-05 14:19:24.019749"
$ws1.Range("F64").Value = "2021-10-05 14:19:24.019752"
$ws1.Range("F65").Value = "2021-10-05 14:19:24.019754"
$ws1.Range("F66").Value = "2021-10-05 14:19:24.019758"
$ws1.Range("F67").Value = "2021-10-05 14:19:24.019761"
$ws1.Range("F68").Value = "2021-10-05 14:19:24.019763"
$ws1.Range("F69").Value = "2021-10-05 14:19:24.019766"
$ws1.Range("F70").Value = "2021-10-05 14:19:24.019768"
$ws1.Range("F71").Value = "2021-10-05 14:19:24.019771"
$ws1.Range("F72").Value = "2021-10-05 14:19:24.019773"
$ws1.Range("F73").Value = "2021-10-05 14:19:24.019776"
$ws1.Range("F74").Value = "2021-10-05 14:19:24.019778"
$ws1.Range("F75").Value = "2021-10-05 14:19:24.019781"
$ws1.Range("F76").Value = "2021-10-05 14:19:24.019783"
$ws1.Range("F77").Value = "2021-10-05 14:19:24.019786"
$ws1.Range("F78").Value = "2021-10-05 14:19:24.019790"
$ws1.Range("F79").Value = "2021-10-05 14:19:24.019793"
$ws1.Range("F80").Value = "2021-10-05 14:19:24.019796"
$ws1.Range("F81").Value = "2021-10-05 14:19:24.019798"
$ws1.Range("F82").Value = "2021-10-05 14:19:24.019801"
$ws1.Range("F83").Value = "2021-10-05 14:19:24.019803"
$ws1.Range("F84").Value = "2021-10-05 14:19:24.019806"
$ws1.Range("F85").Value = "2021-10-05 14:19:24.019808"
$ws1.Range("F86").Value = "2021-10-05 14:19:24.019811"
$ws1.Range("F87").Value = "2021-10-05 14:19:24.019813"
$ws1.Range("F88").Value = "2021-10-05 14:19:24.019816"
$ws1.Range("F89").Value = "2021-10-05 14:19:24.019818"
$ws1.Range("F90").Value = "2021-10-05 14:19:24.019821"
$ws1.Range("F91").Value = "2021-10-05 14:19:24.019823"
$ws1.Range("F92").Value = "2021-10-05 14:19:24.019826"
$ws1.Range("F93").Value = "2021-10-05 14:19:24.019828"
$ws1.Range("F94").Value = "2021-10-05 14:19:24.019832"
$ws1.Range("F95").Value = "2021-10-05 14:19:24.019835"
$ws1.Range("F96").Value = "2021-10-05 14:19:24.019838"
$ws1.Range("F97").Value = "2021-10-05 14:19:24.019840"
$ws1.Range("F98").Value = "2021-10-05 14:19:24.019843"

# Add the new "metadata" worksheet, placed right after "data"
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "metadata"

# Header row (row 1) - bold, thin border, centered, top-aligned to match the
# "data" sheet's header styling
$headerRange = $ws2.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.Item(1).LineStyle = 1
$headerRange.Borders.Item(2).LineStyle = 1
$headerRange.Borders.Item(3).LineStyle = 1
$headerRange.Borders.Item(4).LineStyle = 1
$headerRange.Borders.Item(1).Weight = 2
$headerRange.Borders.Item(2).Weight = 2
$headerRange.Borders.Item(3).Weight = 2
$headerRange.Borders.Item(4).Weight = 2

$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# A2 gets the same header-ish styling as column A in the "data" sheet (index cell)
$a2 = $ws2.Range("A2")
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.Item(1).LineStyle = 1
$a2.Borders.Item(2).LineStyle = 1
$a2.Borders.Item(3).LineStyle = 1
$a2.Borders.Item(4).LineStyle = 1
$a2.Borders.Item(1).Weight = 2
$a2.Borders.Item(2).Weight = 2
$a2.Borders.Item(3).Weight = 2
$a2.Borders.Item(4).Weight = 2
$a2.Value = 0

$ws2.Range("B2").Value = "Cardiac arrhythmias"
$ws2.Range("C2").Value = 842

# D2 must stay the literal text "6.100" rather than become the number 6.1
$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "6.100"

$ws2.Range("E2").Value = "2021-09-28T09:49:50.055065Z"
$ws2.Range("F2").Value = "2021-10-05 14:19:24.016149"
$ws2.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/842/?format=json"

Write-Output "metadata sheet added and timestamps refreshed"
